# Apply corrections to the "currentshop" overlay sheet and add the new
# "Pepper Lights" overlay entry (shifting existing rows down by one).
#
# Resulting table (rows 2-18, columns A:D = Item, Cost, Slot, Priority):
#   2  Pepper Lights     60   2   7
#   3  Coffee Stand      65   9   7
#   4  Newspaper Stand   80   5   7
#   5  Hypno Clock       80  16   7
#   6  Gumball Machine  120   7   7
#   7  Extra Burner     150  14   7
#   8  Ceiling Fan      150   2   7
#   9  TV               150   4   7
#  10  Extra Burner 2   200  15   7
#  11  Arcade Cabinet   400   6   7
#  12  Jukebox          500   8   7
#  13  Royal Crown     1000  17   7
#  14  Doorbell          30   1   8
#  15  Beef Alarm        90   9   8
#  16  Chicken Alarm     90  10   8
#  17  Pork Alarm        90  12   8
#  18  Steak Alarm       90  13   8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currentshop")

$newData = @(
    @("Pepper Lights",    60,  2, 7),
    @("Coffee Stand",     65,  9, 7),
    @("Newspaper Stand",  80,  5, 7),
    @("Hypno Clock",      80, 16, 7),
    @("Gumball Machine", 120,  7, 7),
    @("Extra Burner",    150, 14, 7),
    @("Ceiling Fan",     150,  2, 7),
    @("TV",              150,  4, 7),
    @("Extra Burner 2",  200, 15, 7),
    @("Arcade Cabinet",  400,  6, 7),
    @("Jukebox",         500,  8, 7),
    @("Royal Crown",    1000, 17, 7),
    @("Doorbell",         30,  1, 8),
    @("Beef Alarm",       90,  9, 8),
    @("Chicken Alarm",    90, 10, 8),
    @("Pork Alarm",       90, 12, 8),
    @("Steak Alarm",      90, 13, 8)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = $i + 2
    $row = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
